$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.061.46"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "2.354.30"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "0.682"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").Value = "239.93"
$ws.Range("E6").Value = "  +1.67%  "
$ws.Range("D7").Value = "74.76"
$ws.Range("E7").Value = "  +1.91%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.601"
$ws.Range("E9").Value = "  +6.41%  "
$ws.Range("E10").Value = "  +2.42%  "
$ws.Range("D11").Value = "57.17"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "32.31"
$ws.Range("E12").Value = "  +14.88%  "
$ws.Range("D13").Value = "7.29"
$ws.Range("E13").Value = "  +7.02%  "
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").Value = "2.707.50"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("D18").Value = "2.354.53"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "44.021.12"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("D21").Value = "6.74"
$ws.Range("E21").Value = "  +5.14%  "
$ws.Range("E22").Value = "  -1.57%  "
$ws.Range("D23").Value = "257.28"
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "1.87"
$ws.Range("E25").Value = "  +18.51%  "
$ws.Range("E26").Value = "  -1.35%  "
$ws.Range("D27").Value = "2.51"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "10.76"
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.90"
$ws.Range("E29").Value = "  +1.70%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.24"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("D31").Value = "175.09"
$ws.Range("E31").Value = "  +1.38%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "0.138"
$ws.Range("E32").Value = "  +3.88%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "0.128"
$ws.Range("E33").Value = "  -2.10%  "
$ws.Range("E34").Value = "  +5.63%  "
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("D36").Value = "5.37"
$ws.Range("E36").Value = "  +2.93%  "
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("D38").Value = "2.38"
$ws.Range("E38").Value = "  -2.95%  "
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("E40").Value = "  +4.08%  "
$ws.Range("E41").Value = "  +11.99%  "
$ws.Range("E42").Value = "  +11.77%  "
$ws.Range("D43").Value = "19.33"
$ws.Range("E43").Value = "  +0.76%  "
$ws.Range("D44").Value = "9.07"
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("E46").Value = "  +7.45%  "
$ws.Range("D47").Value = "2.51"
$ws.Range("E47").Value = "  +8.42%  "
$ws.Range("E48").Value = "  +2.62%  "
$ws.Range("D49").Value = "101.12"
$ws.Range("E49").Value = "  +2.83%  "
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("D51").Value = "56.83"
$ws.Range("E51").Value = "  +8.76%  "
